$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 13).Value = 179
$ws.Cells.Item(33, 9).Value = 50
$ws.Cells.Item(33, 11).Value = 50
$ws.Cells.Item(33, 8).Value = 1016.6667
$ws.Cells.Item(41, 13).Value = 165.5
$ws.Cells.Item(41, 14).Value = -1179.75
$ws.Cells.Item(41, 10).Value = 299.75
$ws.Cells.Item(41, 12).Value = 299.75
$ws.Cells.Item(41, 9).Value = 274.5
$ws.Cells.Item(41, 11).Value = 274.5
$ws.Cells.Item(41, 8).Value = 287.125
$ws.Cells.Item(48, 13).ClearContents()
$ws.Cells.Item(48, 14).Value = -24573.5
$ws.Cells.Item(48, 10).Value = 7996.5
$ws.Cells.Item(48, 12).Value = 23989.5
$ws.Cells.Item(48, 9).Value = 0
$ws.Cells.Item(48, 11).Value = 0
$ws.Cells.Item(48, 8).Value = 7996.5
$ws.Cells.Item(56, 13).ClearContents()
$ws.Cells.Item(56, 14).Value = -25057.5
$ws.Cells.Item(56, 10).Value = 7996.5
$ws.Cells.Item(56, 12).Value = 23989.5
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 8).Value = 7996.5
$ws.Cells.Item(70, 13).Value = -10230
$ws.Cells.Item(70, 14).Value = -15165
$ws.Cells.Item(70, 10).Value = 4875
$ws.Cells.Item(70, 12).Value = 14625
$ws.Cells.Item(70, 9).Value = 3500
$ws.Cells.Item(70, 11).Value = 10500
$ws.Cells.Item(70, 8).Value = 4722.222
$ws.Cells.Item(73, 13).Value = -9564
$ws.Cells.Item(73, 14).Value = -16497
$ws.Cells.Item(73, 10).Value = 4875
$ws.Cells.Item(73, 12).Value = 14625
$ws.Cells.Item(73, 9).Value = 3500
$ws.Cells.Item(73, 11).Value = 10500
$ws.Cells.Item(73, 8).Value = 4722.222
$ws.Cells.Item(88, 13).Value = -1677.3333
$ws.Cells.Item(88, 14).Value = -99700
$ws.Cells.Item(88, 10).Value = 98888
$ws.Cells.Item(88, 12).Value = 98888
$ws.Cells.Item(88, 9).Value = 2083.3333
$ws.Cells.Item(88, 11).Value = 2083.3333
$ws.Cells.Item(88, 8).Value = 26284.5
$ws.Cells.Item(91, 13).Value = -679.3332999999998
$ws.Cells.Item(91, 14).Value = -101696
$ws.Cells.Item(91, 10).Value = 98888
$ws.Cells.Item(91, 12).Value = 98888
$ws.Cells.Item(91, 9).Value = 2083.3333
$ws.Cells.Item(91, 11).Value = 2083.3333
$ws.Cells.Item(91, 8).Value = 26284.5
$ws.Cells.Item(107, 13).Value = 1026.6
$ws.Cells.Item(107, 14).Value = -4801.75
$ws.Cells.Item(107, 10).Value = 961.75
$ws.Cells.Item(107, 12).Value = 961.75
$ws.Cells.Item(107, 9).Value = 893.4
$ws.Cells.Item(107, 11).Value = 893.4
$ws.Cells.Item(107, 8).Value = 923.7778
$ws.Cells.Item(132, 13).Value = -3807.0905
$ws.Cells.Item(132, 9).Value = 2112.3635
$ws.Cells.Item(132, 11).Value = 6337.0905
$ws.Cells.Item(132, 8).Value = 2186.3333
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 13).Value = -91.60000000000002
$ws.Cells.Item(74, 9).Value = 965.6
$ws.Cells.Item(74, 11).Value = 965.6
$ws.Cells.Item(74, 8).Value = 963
$ws.Cells.Item(77, 13).Value = -460
$ws.Cells.Item(77, 9).Value = 965.6
$ws.Cells.Item(77, 11).Value = 4828
$ws.Cells.Item(77, 8).Value = 963
$ws.Cells.Item(88, 13).Value = 216
$ws.Cells.Item(88, 9).Value = 190
$ws.Cells.Item(88, 11).Value = 190
$ws.Cells.Item(88, 8).Value = 2891.2
$ws.Cells.Item(91, 13).Value = 1214
$ws.Cells.Item(91, 9).Value = 190
$ws.Cells.Item(91, 11).Value = 190
$ws.Cells.Item(91, 8).Value = 2891.2
$ws.Cells.Item(132, 13).Value = -6464
$ws.Cells.Item(132, 9).Value = 2998
$ws.Cells.Item(132, 11).Value = 8994
$ws.Cells.Item(132, 8).Value = 2998
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 13).Value = -3215.4
$ws.Cells.Item(86, 9).Value = 4338.4
$ws.Cells.Item(86, 11).Value = 4338.4
$ws.Cells.Item(86, 8).Value = 5466.5557
$ws.Cells.Item(89, 13).Value = -16076
$ws.Cells.Item(89, 9).Value = 4338.4
$ws.Cells.Item(89, 11).Value = 21692
$ws.Cells.Item(89, 8).Value = 5466.5557
$ws.Cells.Item(134, 13).Value = -9402
$ws.Cells.Item(134, 14).Value = -44067
$ws.Cells.Item(134, 10).Value = 12999
$ws.Cells.Item(134, 12).Value = 38997
$ws.Cells.Item(134, 9).Value = 3979
$ws.Cells.Item(134, 11).Value = 11937
$ws.Cells.Item(134, 8).Value = 4730.6665
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(29, 14).ClearContents()
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(29, 8).Value = 0
$ws.Cells.Item(31, 13).Value = -1685.2727
$ws.Cells.Item(31, 9).Value = 1980.2727
$ws.Cells.Item(31, 11).Value = 1980.2727
$ws.Cells.Item(31, 8).Value = 2307.2666
$ws.Cells.Item(34, 13).Value = -1778.2727
$ws.Cells.Item(34, 9).Value = 1980.2727
$ws.Cells.Item(34, 11).Value = 1980.2727
$ws.Cells.Item(34, 8).Value = 2307.2666
$ws.Cells.Item(68, 14).Value = -41489.43
$ws.Cells.Item(68, 10).Value = 39991.43
$ws.Cells.Item(68, 12).Value = 39991.43
$ws.Cells.Item(68, 8).Value = 39991.43
$ws.Cells.Item(71, 14).Value = -127462.29
$ws.Cells.Item(71, 10).Value = 39991.43
$ws.Cells.Item(71, 12).Value = 119974.29
$ws.Cells.Item(71, 8).Value = 39991.43
$ws.Cells.Item(99, 13).Value = -3121.8
$ws.Cells.Item(99, 9).Value = 4619.8
$ws.Cells.Item(99, 11).Value = 4619.8
$ws.Cells.Item(99, 8).Value = 5161.5557
$ws.Cells.Item(102, 14).Value = -54868
$ws.Cells.Item(102, 10).Value = 50000
$ws.Cells.Item(102, 12).Value = 50000
$ws.Cells.Item(102, 8).Value = 50000
$ws.Cells.Item(126, 13).Value = -11389.4
$ws.Cells.Item(126, 9).Value = 4619.8
$ws.Cells.Item(126, 11).Value = 13859.4
$ws.Cells.Item(126, 8).Value = 5161.5557
$ws.Cells.Item(134, 13).Value = -1965
$ws.Cells.Item(134, 9).Value = 1500
$ws.Cells.Item(134, 11).Value = 4500
$ws.Cells.Item(134, 8).Value = 1721
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(128, 13).Value = -969700.98
$ws.Cells.Item(128, 9).Value = 324893.66
$ws.Cells.Item(128, 11).Value = 974680.98
$ws.Cells.Item(128, 8).Value = 324893.66
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 13).Value = -333333070
$ws.Cells.Item(70, 9).Value = 333333340
$ws.Cells.Item(70, 11).Value = 333333340
$ws.Cells.Item(70, 8).Value = 250001890
$ws.Cells.Item(73, 13).Value = -333332404
$ws.Cells.Item(73, 9).Value = 333333340
$ws.Cells.Item(73, 11).Value = 333333340
$ws.Cells.Item(73, 8).Value = 250001890
$ws.Cells.Item(122, 13).Value = -2778.5002
$ws.Cells.Item(122, 9).Value = 1742.8334
$ws.Cells.Item(122, 11).Value = 5228.5002
$ws.Cells.Item(122, 8).Value = 1775.1428
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 13).Value = -8059.999899999999
$ws.Cells.Item(122, 14).Value = -18400
$ws.Cells.Item(122, 10).Value = 4500
$ws.Cells.Item(122, 12).Value = 13500
$ws.Cells.Item(122, 9).Value = 3503.3333
$ws.Cells.Item(122, 11).Value = 10509.9999
$ws.Cells.Item(122, 8).Value = 3645.7144
$ws.Cells.Item(132, 13).Value = -19015.625
$ws.Cells.Item(132, 9).Value = 7181.875
$ws.Cells.Item(132, 11).Value = 21545.625
$ws.Cells.Item(132, 8).Value = 8823.200000000001
$ws.Cells.Item(136, 13).Value = -3036.529500000001
$ws.Cells.Item(136, 9).Value = 1862.1765
$ws.Cells.Item(136, 11).Value = 5586.529500000001
$ws.Cells.Item(136, 8).Value = 2036.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 14).ClearContents()
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 8).Value = 9999999
$ws.Cells.Item(81, 13).Value = -2215.75
$ws.Cells.Item(81, 14).Value = -2560
$ws.Cells.Item(81, 10).Value = 219
$ws.Cells.Item(81, 12).Value = 438
$ws.Cells.Item(81, 9).Value = 1638.375
$ws.Cells.Item(81, 11).Value = 3276.75
$ws.Cells.Item(81, 8).Value = 1480.6666
$ws.Cells.Item(84, 13).Value = -11079.75
$ws.Cells.Item(84, 14).Value = -12798
$ws.Cells.Item(84, 10).Value = 219
$ws.Cells.Item(84, 12).Value = 2190
$ws.Cells.Item(84, 9).Value = 1638.375
$ws.Cells.Item(84, 11).Value = 16383.75
$ws.Cells.Item(84, 8).Value = 1480.6666
$ws.Cells.Item(122, 13).Value = -8596
$ws.Cells.Item(122, 14).Value = -18397.75
$ws.Cells.Item(122, 10).Value = 4499.25
$ws.Cells.Item(122, 12).Value = 13497.75
$ws.Cells.Item(122, 9).Value = 3682
$ws.Cells.Item(122, 11).Value = 11046
$ws.Cells.Item(122, 8).Value = 4008.9

Write-Host "Applied 191 cell updates across 8 sheets."